# Auto-generated Excel COM-interop script to apply the cryptos.xlsx cell updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) contain values that look numeric (e.g. '1.000',
# '246.07') but must be preserved as literal text, exactly like the source data.
# Pre-format the ranges as Text so assigning the string values below doesn't get
# silently coerced into real numbers, then restore the default style afterwards so
# no stray formatting is left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '31.115.58', '  +1.77%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.953.22', '  +0.78%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  +0.08%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '246.07', '  -0.22%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  +0.11%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4875', '  +0.54%  ')
    ,@(8, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '44.59', '  -0.09%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2964', '  +1.19%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06819', '  -0.01%  ')
    ,@(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.04', '  -2.39%  ')
    ,@(12, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '105.88', '  -6.35%  ')
    ,@(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.935.79', '  -0.04%  ')
    ,@(14, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07723', '  +1.66%  ')
    ,@(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.398', '  -1.93%  ')
    ,@(16, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7115', '  +4.22%  ')
    ,@(17, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '284.50', '  -4.81%  ')
    ,@(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '31.003.32', '  +1.38%  ')
    ,@(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007736', '  +0.79%  ')
    ,@(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.20', '  +0.42%  ')
    ,@(21, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.192.70', '  +0.23%  ')
    ,@(22, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9998', '  +0.09%  ')
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.520', '  -1.05%  ')
    ,@(24, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.14%  ')
    ,@(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.596', '  +0.98%  ')
    ,@(26, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.899', '  +3.43%  ')
    ,@(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '168.37', '  +0.14%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.87', '  -2.89%  ')
    ,@(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.188', '  +2.08%  ')
    ,@(30, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1051', '  -1.93%  ')
    ,@(31, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.440', '  +0.34%  ')
    ,@(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.723', '  +14.66%  ')
    ,@(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.466', '  +6.71%  ')
    ,@(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04993', '  -0.75%  ')
    ,@(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7599', '  +1.34%  ')
    ,@(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.161', '  +0.66%  ')
    ,@(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.731', '  +0.78%  ')
    ,@(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02039', '  -0.16%  ')
    ,@(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.714', '  +0.49%  ')
    ,@(40, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.143', '  +4.95%  ')
    ,@(41, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.419', '  +9.45%  ')
    ,@(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4471', '  -0.04%  ')
    ,@(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '109.57', '  -0.71%  ')
    ,@(44, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8813', '  +0.78%  ')
    ,@(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '72.68', '  +3.95%  ')
    ,@(46, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9992', '  -0.15%  ')
    ,@(47, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.467', '  +2.02%  ')
    ,@(48, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '982.77', '  +15.98%  ')
    ,@(49, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1276', '  +3.13%  ')
    ,@(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.385', '  +0.95%  ')
    ,@(51, 'WOONetwork', 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo', '0.2588', '  +1.41%  ')
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
}

# Restore the default ('Normal') cell style now that the text values are locked in,
# so the worksheet's formatting matches the unmodified original.
$priceRange.Style = "Normal"
$volRange.Style = "Normal"

